# Region mapping normalization: rename the two header labels in Sheet1
# ("Region starting 2016" -> "Region until 2016",
#  "Region until 2017"    -> "Region after 2016"),
# matching the updated ETL pipeline's region-normalization logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B1").Value = "Region until 2016"
$ws.Range("C1").Value = "Region after 2016"

# Leave the cursor where the author ended up after editing the headers.
$ws.Range("F6").Select()
